# run prepare & render with final data
# Updates the computed statistics (mean / CI_low / CI_high) on Sheet 1
# with refreshed values from the final model run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0775236584228791
$ws.Range("C2").Value = 0.0587637933184519
$ws.Range("D2").Value = 0.0962835235273064

$ws.Range("B3").Value = 0.0825944663024482
$ws.Range("C3").Value = 0.0547133371389904
$ws.Range("D3").Value = 0.110475595465906

$ws.Range("B4").Value = 0.0495344129192724
$ws.Range("C4").Value = -0.0118154213608932

$ws.Range("C5").Value = -0.00960239297963123

$ws.Range("C6").Value = 0.0728607170371897
$ws.Range("D6").Value = 0.215814026726263

$ws.Range("C7").Value = 0.0293215892202031

$ws.Range("B8").Value = 0.139448173442151
$ws.Range("C8").Value = 0.0580077114529696
$ws.Range("D8").Value = 0.220888635431333

$ws.Range("B9").Value = 0.0457932964454881

$ws.Range("B10").Value = 0.039019983833404
$ws.Range("C10").Value = -0.0505289947955081

$ws.Range("B11").Value = 0.0174473111262245
$ws.Range("C11").Value = -0.0152238153386416
$ws.Range("D11").Value = 0.0501184375910905

$ws.Range("B12").Value = 0.0793351577642075
$ws.Range("C12").Value = 0.0215888870035115
$ws.Range("D12").Value = 0.137081428524904

$ws.Range("B13").Value = 0.0572671934345981

$ws.Range("B14").Value = 0.100905525447827
$ws.Range("C14").Value = 0.064886266825133
$ws.Range("D14").Value = 0.136924784070522
